$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 118
$ws1.Range("F3").Value = 2136
$ws1.Range("F4").Value = 21
$ws1.Range("F5").Value = 11108
$ws1.Range("F10").Value = 11007
$ws1.Range("F12").Value = 1137
$ws1.Range("F13").Value = 35
$ws1.Range("F15").Value = 5513
$ws1.Range("F17").Value = 3426

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 118
$ws4.Range("F3").Value = 2136
$ws4.Range("F5").Value = 21
$ws4.Range("F7").Value = 11108
$ws4.Range("F12").Value = 11007
$ws4.Range("F14").Value = 1137
$ws4.Range("F15").Value = 35
$ws4.Range("F17").Value = 5513
$ws4.Range("F19").Value = 3426
